# Applies the Phoenix_Profits market-data refresh described in the commit diff.
# Values are plain numeric cell updates (no formulas); one cell (CUL!M122) is
# cleared entirely because the refreshed row no longer has a value there.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 91370.27
$ws.Range("I2").Value = 213.28572
$ws.Range("J2").Value = 250895
$ws.Range("K2").Value = 213.28572
$ws.Range("L2").Value = 250895
$ws.Range("M2").Value = -100.28572
$ws.Range("N2").Value = -251121
# Row 28
$ws.Range("H28").Value = 1608.5
$ws.Range("I28").Value = 1917.125
$ws.Range("J28").Value = 991.25
$ws.Range("K28").Value = 1917.125
$ws.Range("L28").Value = 991.25
$ws.Range("M28").Value = -1432.125
$ws.Range("N28").Value = -1961.25
# Row 76
$ws.Range("H76").Value = 4374
$ws.Range("I76").Value = 4374
$ws.Range("K76").Value = 4374
$ws.Range("M76").Value = -4059
# Row 79
$ws.Range("H79").Value = 4374
$ws.Range("I79").Value = 4374
$ws.Range("K79").Value = 4374
$ws.Range("M79").Value = -3282
# Row 80
$ws.Range("H80").Value = 335.7
$ws.Range("J80").Value = 997.5
$ws.Range("L80").Value = 2992.5
$ws.Range("N80").Value = -4988.5
# Row 83
$ws.Range("H83").Value = 335.7
$ws.Range("J83").Value = 997.5
$ws.Range("L83").Value = 8977.5
$ws.Range("N83").Value = -18961.5
# Row 86
$ws.Range("H86").Value = 3854.1428
$ws.Range("I86").Value = 2995
$ws.Range("K86").Value = 2995
$ws.Range("M86").Value = -1872
# Row 88
$ws.Range("H88").Value = 1385
$ws.Range("J88").Value = 1385
$ws.Range("L88").Value = 1385
$ws.Range("N88").Value = -2197
# Row 89
$ws.Range("H89").Value = 3854.1428
$ws.Range("I89").Value = 2995
$ws.Range("K89").Value = 14975
$ws.Range("M89").Value = -9359
# Row 91
$ws.Range("H91").Value = 1385
$ws.Range("J91").Value = 1385
$ws.Range("L91").Value = 1385
$ws.Range("N91").Value = -4193
# Row 111
$ws.Range("H111").Value = 55557836
$ws.Range("J111").Value = 166668540
$ws.Range("L111").Value = 500005620
$ws.Range("N111").Value = -500011754
# Row 138
$ws.Range("H138").Value = 2250.3333
$ws.Range("I138").Value = 1289.8
$ws.Range("K138").Value = 3869.4
$ws.Range("M138").Value = 1270.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 303942.66
$ws.Range("I74").Value = 420415.53
$ws.Range("J74").Value = 70996.92
$ws.Range("K74").Value = 420415.53
$ws.Range("L74").Value = 70996.92
$ws.Range("M74").Value = -419541.53
$ws.Range("N74").Value = -72744.92
# Row 77
$ws.Range("H77").Value = 303942.66
$ws.Range("I77").Value = 420415.53
$ws.Range("J77").Value = 70996.92
$ws.Range("K77").Value = 2102077.65
$ws.Range("L77").Value = 354984.6
$ws.Range("M77").Value = -2097709.65
$ws.Range("N77").Value = -363720.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 27029494
$ws.Range("I86").Value = 34485116
$ws.Range("J86").Value = 2862.125
$ws.Range("K86").Value = 34485116
$ws.Range("L86").Value = 2862.125
$ws.Range("M86").Value = -34483993
$ws.Range("N86").Value = -5108.125
# Row 89
$ws.Range("H89").Value = 27029494
$ws.Range("I89").Value = 34485116
$ws.Range("J89").Value = 2862.125
$ws.Range("K89").Value = 172425580
$ws.Range("L89").Value = 14310.625
$ws.Range("M89").Value = -172419964
$ws.Range("N89").Value = -25542.625
# Row 96
$ws.Range("H96").Value = 27000
$ws.Range("J96").Value = 44000
$ws.Range("L96").Value = 44000
$ws.Range("N96").Value = -49492
# Row 134
$ws.Range("H134").Value = 2109.8032
$ws.Range("I134").Value = 1814.0667
$ws.Range("K134").Value = 5442.2001
$ws.Range("M134").Value = -2907.2001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1373.3793
$ws.Range("I31").Value = 1335.3405
$ws.Range("J31").Value = 1535.909
$ws.Range("K31").Value = 1335.3405
$ws.Range("L31").Value = 1535.909
$ws.Range("M31").Value = -1040.3405
$ws.Range("N31").Value = -2125.909
# Row 34
$ws.Range("H34").Value = 1373.3793
$ws.Range("I34").Value = 1335.3405
$ws.Range("J34").Value = 1535.909
$ws.Range("K34").Value = 1335.3405
$ws.Range("L34").Value = 1535.909
$ws.Range("M34").Value = -1133.3405
$ws.Range("N34").Value = -1939.909
# Row 134
$ws.Range("H134").Value = 2898.4048
$ws.Range("I134").Value = 3614.8845
$ws.Range("K134").Value = 10844.6535
$ws.Range("M134").Value = -8309.6535

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 365.14285
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 365.14285
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3286.28565
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -8186.28565

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3500
$ws.Range("I80").Value = 3500
$ws.Range("K80").Value = 3500
$ws.Range("M80").Value = -2502
# Row 83
$ws.Range("H83").Value = 3500
$ws.Range("I83").Value = 3500
$ws.Range("K83").Value = 17500
$ws.Range("M83").Value = -12508
# Row 102
$ws.Range("H102").Value = 29707.166
$ws.Range("I102").Value = 40055.64
$ws.Range("J102").Value = 6187.909
$ws.Range("K102").Value = 40055.64
$ws.Range("L102").Value = 6187.909
$ws.Range("M102").Value = -38433.64
$ws.Range("N102").Value = -9431.909

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 20230.39
$ws.Range("J7").Value = 3423.5715
$ws.Range("L7").Value = 3423.5715
$ws.Range("N7").Value = -3647.5715
# Row 82
$ws.Range("H82").Value = 1821
$ws.Range("I82").Value = 1789.1818
$ws.Range("J82").Value = 1996
$ws.Range("K82").Value = 1789.1818
$ws.Range("L82").Value = 1996
$ws.Range("M82").Value = -1428.1818
$ws.Range("N82").Value = -2718
# Row 85
$ws.Range("H85").Value = 1821
$ws.Range("I85").Value = 1789.1818
$ws.Range("J85").Value = 1996
$ws.Range("K85").Value = 1789.1818
$ws.Range("L85").Value = 1996
$ws.Range("M85").Value = -541.1818000000001
$ws.Range("N85").Value = -4492
# Row 126
$ws.Range("H126").Value = 20230.39
$ws.Range("J126").Value = 3423.5715
$ws.Range("L126").Value = 10270.7145
$ws.Range("N126").Value = -15210.7145
# Row 132
$ws.Range("H132").Value = 3177.878
$ws.Range("I132").Value = 2551.862
$ws.Range("J132").Value = 4690.75
$ws.Range("K132").Value = 7655.586
$ws.Range("L132").Value = 14072.25
$ws.Range("M132").Value = -5125.586
$ws.Range("N132").Value = -19132.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 34487356
$ws.Range("I126").Value = 38465876
$ws.Range("K126").Value = 115397628
$ws.Range("M126").Value = -115395158

